# Final commit of upload excel file:
# update a handful of contact-detail values and restore the row height
# that Excel applies when the sheet is re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FirstName for row 2 was blank, now filled in
$ws.Range("B2").Value = "rohan"

# Email for row 2 corrected
$ws.Range("J2").Value = "tintu@gmail.com"

# Hobbies for row 2: dropped "Driving"
$ws.Range("L2").Value = "Reading ,Drawing"

# FirstName for row 3 changed
$ws.Range("B3").Value = "mini"

# Street for row 3 changed
$ws.Range("H3").Value = "abcd"

# Hobbies for row 3: trailing comma removed
$ws.Range("L3").Value = "Reading ,Writing"

# Row height for the header + both data rows bumped slightly
$ws.Rows("1:3").RowHeight = 19.5
